# Add data for 2022-07-23
# - Rename sheet to reflect new "through" date
# - Update the "2022 (through ...)" column header label
# - Update the current-year (column I) figures for July (row 7) and
#   August (row 8), plus the Total row (I14), to match the new sum.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "Through 2022-07-15"

# Update the header label in column I, row 1
$ws.Range("I1").Value = "2022 (through 07-15)"

# Update the monthly figures that changed
$ws.Range("I7").Value = 142
$ws.Range("I8").Value = 82

# Update the Total row to reflect the new sum
$ws.Range("I14").Value = 887
